# This edit re-shuffles the weekly "Fruta / hortaliza" price-report rows
# (rows 2-30 of the only worksheet): the per-observation data in columns
# D,L,M,N,O,P,Q,R,S,T of each row is replaced by the data that used to live
# in a different row (a like-for-like permutation of the weekly
# observations), while columns A,B,C,E,F,G,H,I,J,K (market/product identity,
# which is identical for every row) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: destination row -> source row whose values it should adopt.
$rowSource = @{}
$rowSource[2] = 12
$rowSource[3] = 20
$rowSource[4] = 29
$rowSource[5] = 3
$rowSource[6] = 19
$rowSource[7] = 25
$rowSource[8] = 4
$rowSource[9] = 8
$rowSource[10] = 9
$rowSource[11] = 17
$rowSource[12] = 22
$rowSource[13] = 23
$rowSource[14] = 26
$rowSource[15] = 2
$rowSource[16] = 6
$rowSource[17] = 7
$rowSource[18] = 24
$rowSource[19] = 16
$rowSource[20] = 5
$rowSource[21] = 15
$rowSource[22] = 10
$rowSource[23] = 11
$rowSource[24] = 30
$rowSource[25] = 13
$rowSource[26] = 14
$rowSource[27] = 18
$rowSource[28] = 27
$rowSource[29] = 28
$rowSource[30] = 21

# Columns holding the per-observation data (Fecha, Calidad, Volumen,
# Precio minimo/maximo/promedio, Unidad de comercializacion, Origen,
# Precio $/Kg, Kg / unidad) that gets permuted across rows.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# 1) Snapshot every current ("before") value first. Several rows feed into
#    each other (the permutation has multi-row cycles), so all reads must
#    happen before any writes to avoid clobbering values still needed.
$snapshot = @{}
foreach ($row in 2..30) {
    foreach ($col in $cols) {
        $ref = "$col$row"
        $snapshot[$ref] = $ws.Range($ref).Value2
    }
}

# 2) Write the permuted values into place.
foreach ($row in 2..30) {
    $src = $rowSource[$row]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $snapshot["$col$src"]
    }
}
